$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.559.86'
$ws.Range('E2').Value = '  +2.34%  '

$ws.Range('D3').Value = '1.638.01'
$ws.Range('E3').Value = '  +3.88%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9962'
$ws.Range('E4').Value = '  -0.56%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.71'
$ws.Range('E5').Value = '  +2.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9962'
$ws.Range('E6').Value = '  -0.75%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3788'
$ws.Range('E7').Value = '  +1.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '53.11'
$ws.Range('E8').Value = '  +6.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3679'
$ws.Range('E9').Value = '  +3.08%  '

$ws.Range('E10').Value = '  +6.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08204'
$ws.Range('E11').Value = '  +3.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9960'
$ws.Range('E12').Value = '  -0.57%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.39'
$ws.Range('E13').Value = '  +7.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.679'
$ws.Range('E14').Value = '  +3.76%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001289'
$ws.Range('E15').Value = '  +5.92%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.488'
$ws.Range('E16').Value = '  +2.95%  '

$ws.Range('D17').Value = '1.634.76'
$ws.Range('E17').Value = '  +3.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.04'
$ws.Range('E18').Value = '  +3.68%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06951'
$ws.Range('E19').Value = '  +3.01%  '

$ws.Range('E20').Value = '  +4.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.619'
$ws.Range('E21').Value = '  +4.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9959'
$ws.Range('E22').Value = '  -0.73%  '

$ws.Range('D23').Value = '23.568.11'
$ws.Range('E23').Value = '  +2.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.02'
$ws.Range('E24').Value = '  +2.82%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.155'
$ws.Range('E25').Value = '  +12.96%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.427'
$ws.Range('E26').Value = '  +2.65%  '

$ws.Range('E27').Value = '  +4.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.37'
$ws.Range('E28').Value = '  +2.91%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.330'
$ws.Range('E29').Value = '  +3.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '136.72'
$ws.Range('E30').Value = '  +3.97%  '

$ws.Range('E31').Value = '  +4.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.898'
$ws.Range('E32').Value = '  +6.76%  '

$ws.Range('D33').Value = '1.813.67'
$ws.Range('E33').Value = '  +3.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9805'
$ws.Range('E34').Value = '  +5.84%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02837'
$ws.Range('E35').Value = '  +6.86%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.48'
$ws.Range('E36').Value = '  +6.60%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.07502'
$ws.Range('E37').Value = '  +2.72%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.246'
$ws.Range('E38').Value = '  +5.11%  '

$ws.Range('E39').Value = '  +2.65%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08870'
$ws.Range('E40').Value = '  +1.63%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.411'
$ws.Range('E41').Value = '  +5.60%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7193'
$ws.Range('E42').Value = '  +5.37%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.77'
$ws.Range('E43').Value = '  +8.64%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.26'
$ws.Range('E44').Value = '  +10.88%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6658'
$ws.Range('E45').Value = '  +5.27%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.374'
$ws.Range('E46').Value = '  +6.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.038'
$ws.Range('E47').Value = '  +1.89%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9953'
$ws.Range('E48').Value = '  -0.71%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08062'
$ws.Range('E49').Value = '  +2.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.59'
$ws.Range('E50').Value = '  +1.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.221'
$ws.Range('E51').Value = '  +3.65%  '
